# Applies the "PID + Complementary Filter + It compiles!" edit:
#  - Collapses several runs that were split apart only for
#    <w:proofErr> spell-check bookmarking back into single runs
#    (the visible text is unchanged, only the run/proofErr
#    structure is cleaned up).
#  - Appends a new 15th October, 2023 journal entry at the end.

$d = $word.ActiveDocument

function Set-ParagraphXml($Marker, $InnerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($Marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Marker not found: $Marker"
    }
    $para = $rng.Paragraphs(1)
    $pr = $para.Range
    $pr.MoveEnd(1, -1)
    $pr.Text = ""
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pr.InsertXML($frag)
}

$langRPr = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'
$pPrLang = '<w:pPr>' + $langRPr + '</w:pPr>'

# 1. "Self Balancing" + " Bot Journal." -> single run
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t>Self Balancing Bot Journal.</w:t></w:r>'
Set-ParagraphXml "Self Balancing" $inner

# 2. "Sending debug info via Bluetooth..." paragraph -> single run
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t>Sending debug info via Bluetooth. I think im gonna create a python application to collect this data so I can save it in files. Debug info is probably going to be very useful. Ill have to see. I have maintained the code neatly, by splitting the code into various headers and translation units. It does compile.</w:t></w:r>'
Set-ParagraphXml "Sending debug info via Bluetooth" $inner

# 3. "I have the MPU6050 already..." -> single run
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t>I have the MPU6050 already. Poor soldering but it works. Lets get testing.</w:t></w:r>'
Set-ParagraphXml "I have the MPU6050 already" $inner

# 4. "As a last minute change..." -> merged first run, trailing "." run kept separate
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t>As a last minute change I might shift from Adafruit_MPU6050 to MPU6050 by ElectronicCats</w:t></w:r>' + '<w:r>' + $langRPr + '<w:t>.</w:t></w:r>'
Set-ParagraphXml "As a last minute change" $inner

# 5. "Im getting much more accurate angle readings..." -> single run
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t>Im getting much more accurate angle readings through accelerometer from there. When the chip is almost upright in either direction, im getting up 89 degrees in both sides. This seems more accurate than what im using right now.</w:t></w:r>'
Set-ParagraphXml "getting much more accurate angle readings" $inner

# 6. "I need to write better debug info..." -> 3 runs
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t xml:space="preserve">I need to write better debug info. I doubt assembling will be done even by tomorrow. </w:t></w:r>' + '<w:r>' + $langRPr + '<w:t xml:space="preserve">Ill have time to make some updates throughout this week. </w:t></w:r>' + '<w:r>' + $langRPr + '<w:t>We are 1 team member short. But its okay. I always felt like the team had only 3 members in the first place.</w:t></w:r>'
Set-ParagraphXml "I need to write better debug info" $inner

# 7. "Now I have the right offsets..." -> single run
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t>Now I have the right offsets which will make my readings accurate. Ill try again, using angular velocity values to find angle.</w:t></w:r>'
Set-ParagraphXml "Now I have the right offsets" $inner

# 8. "... (since stage 2 involves moving the robo)" bullet -> keep pStyle/numPr, merge trailing runs
$kalmanPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
$inner = $kalmanPPr + '<w:r>' + $langRPr + '<w:t>Kalman Filter library to combine gyro + accel to get better angle values.</w:t></w:r>' + '<w:r>' + $langRPr + '<w:t xml:space="preserve"> (since stage 2 involves moving the robo)</w:t></w:r>'
Set-ParagraphXml "Kalman Filter library to combine gyro" $inner

# 9. "The mid term 2 is finished..." -> single run
$inner = $pPrLang + '<w:r>' + $langRPr + '<w:t xml:space="preserve">The mid term 2 is finished. Im back to working on the thing. The above three points are to be worked on. </w:t></w:r>'
Set-ParagraphXml "The mid term 2 is finished" $inner

# 10. Append the new 15th October, 2023 journal entry at the end of the document.
$newParas = '<w:p>' + $pPrLang + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>15</w:t></w:r>' + '<w:r><w:rPr><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr><w:t>th</w:t></w:r>' + '<w:r>' + $langRPr + '<w:t xml:space="preserve"> October, 2023</w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>I am finally done with the code. I settled with – a complementary filter. No Kalman. There is a PID library which I am using to generate control.</w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>My guess is that I do not need integral term for the PID. I need to tune P and D terms only. I will be 0.</w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>But I cannot say anything for sure until we actually try it out.</w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>At the very least, the behaviour of PID is as predicted. As the offset increases, it generates a higher control signal. And it also generates it in the opposite direction.</w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>The sad news is that im too tired right now to work on the interface.</w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t xml:space="preserve">There is no happy news. </w:t></w:r>' + '</w:p>'
$newParas += '<w:p>' + $pPrLang + '<w:r>' + $langRPr + '<w:t>Maybe ill work on the Bluetooth communication thing by tonight.</w:t></w:r>' + '<w:r>' + $langRPr + '<w:t xml:space="preserve"> Lets see.</w:t></w:r>' + '</w:p>'

$endRng = $d.Content
$endRng.Collapse(0)
$endFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParas + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRng.InsertXML($endFrag)

Write-Output "Edit complete."
